$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume cells are stored as plain text in the source data (e.g. "64.507.35",
# "0.460", "26.20", "0.0000228" ...). Pin each touched cell to Text format before
# writing so Excel does not reinterpret/renormalize numeric-looking values (which
# would silently drop trailing zeros, switch to scientific notation, etc.) and the
# cell keeps matching the rest of the (text-only) column.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = '65.005.08'
$ws.Range("E2").Value = '  +6.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = '3.108.75'
$ws.Range("E3").Value = '  +4.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = '557.98'
$ws.Range("E5").Value = '  +2.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = '145.13'
$ws.Range("E6").Value = '  +12.98%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = '3.107.46'
$ws.Range("E8").Value = '  +5.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  +1.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = '7.04'
$ws.Range("E10").Value = '  +19.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.153'
$ws.Range("E11").Value = '  +7.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("E12").Value = '  +5.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +5.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '35.43'
$ws.Range("E14").Value = '  +6.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = '3.581.61'
$ws.Range("E15").Value = '  +3.67%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = '64.840.81'
$ws.Range("E16").Value = '  +6.25%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = '3.090.60'
$ws.Range("E17").Value = '  +3.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = '0.109'
$ws.Range("E18").Value = '  -0.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = '6.75'
$ws.Range("E19").Value = '  +3.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D20").Value = '485.47'
$ws.Range("E20").Value = '  +3.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = '13.74'
$ws.Range("E21").Value = '  +6.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = '0.675'
$ws.Range("E22").Value = '  +3.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = '7.55'
$ws.Range("E23").Value = '  +10.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = '13.22'
$ws.Range("E24").Value = '  +12.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D25").Value = '80.97'
$ws.Range("E25").Value = '  +2.44%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value = '2.82'
$ws.Range("E27").Value = '  +5.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D28").Value = '7.98'
$ws.Range("E28").Value = '  +6.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D29").Value = '2.06'
$ws.Range("E29").Value = '  +10.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D30").Value = '0.996'
$ws.Range("E30").Value = '  -0.68%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D31").Value = '26.20'
$ws.Range("E31").Value = '  +4.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D32").Value = '1.16'
$ws.Range("E32").Value = '  +4.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D33").Value = '2.45'
$ws.Range("E33").Value = '  +8.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D34").Value = '5.75'
$ws.Range("E34").Value = '  +8.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D35").Value = '55.29'
$ws.Range("E35").Value = '  +2.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D36").Value = '6.13'
$ws.Range("E36").Value = '  +6.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D37").Value = '468.54'
$ws.Range("E37").Value = '  +5.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0409'
$ws.Range("E38").Value = '  +9.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0829'
$ws.Range("E39").Value = '  +6.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = '3.030.55'
$ws.Range("E40").Value = '  -2.19%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = '8.31'
$ws.Range("E42").Value = '  +4.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = '2.73'
$ws.Range("E43").Value = '  +20.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = '28.04'
$ws.Range("E44").Value = '  +12.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.260'
$ws.Range("E45").Value = '  +10.09%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D47").Value = '2.09'
$ws.Range("E47").Value = '  +10.17%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.0₃0521'
$ws.Range("E49").Value = '  +11.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '118.07'
$ws.Range("E50").Value = '  +4.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("D51").Value = '2.07'
$ws.Range("E51").Value = '  +6.19%  '
